$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 05:48:36"
$ws.Range("E3").Value = "2026-02-20 05:48:39"
$ws.Range("H3").Value = "'95%"
$ws.Range("I3").Value = "1.1 mm"
$ws.Range("N3").Value = "-6.2 °C 5:02 TU"
$ws.Range("O3").Value = "-6.0 °C"
$ws.Range("E4").Value = "2026-02-20 05:48:41"
$ws.Range("H4").Value = "'54%"
$ws.Range("J4").Value = "1018.5 hPa"
$ws.Range("N4").Value = "5.5 °C 5:11 TU"
$ws.Range("O4").Value = "8.6 °C"
$ws.Range("E5").Value = "2026-02-20 05:48:44"
$ws.Range("N5").Value = "-5.9 °C 5:13 TU"
$ws.Range("E6").Value = "2026-02-20 05:48:47"
$ws.Range("H6").Value = "'81%"
$ws.Range("J6").Value = "1018.5 hPa"
$ws.Range("N6").Value = "3.1 °C 5:05 TU"
$ws.Range("O6").Value = "5.4 °C"
$ws.Range("E7").Value = "2026-02-20 05:48:49"
$ws.Range("H7").Value = "'48%"
$ws.Range("J7").Value = "1018.3 hPa"
$ws.Range("N7").Value = "10.2 °C 5:28 TU"
$ws.Range("O7").Value = "11.0 °C"
$ws.Range("E8").Value = "2026-02-20 05:48:52"
$ws.Range("J8").Value = "1019.2 hPa"
$ws.Range("N8").Value = "6.4 °C 5:18 TU"
$ws.Range("E9").Value = "2026-02-20 05:48:55"
$ws.Range("L9").Value = "50.4 km/h - 357º 5:12 TU"
$ws.Range("E10").Value = "2026-02-20 05:48:57"
$ws.Range("E11").Value = "2026-02-20 05:49:00"
$ws.Range("E12").Value = "2026-02-20 05:49:03"
$ws.Range("H12").Value = "'52%"
$ws.Range("E13").Value = "2026-02-20 05:49:05"
$ws.Range("J13").Value = "1020.0 hPa"
$ws.Range("L13").Value = "62.6 km/h - 34º 5:15 TU"
$ws.Range("N13").Value = "4.3 °C 5:27 TU"
$ws.Range("E14").Value = "2026-02-20 05:49:08"
$ws.Range("E15").Value = "2026-02-20 05:49:11"
$ws.Range("O15").Value = "12.7 °C"
$ws.Range("E16").Value = "2026-02-20 05:49:14"
$ws.Range("H16").Value = "'57%"
$ws.Range("E17").Value = "2026-02-20 05:49:16"
$ws.Range("E18").Value = "2026-02-20 05:49:19"
$ws.Range("J18").Value = "1018.9 hPa"
$ws.Range("N18").Value = "0.4 °C 5:28 TU"
$ws.Range("O18").Value = "1.9 °C"
$ws.Range("E19").Value = "2026-02-20 05:49:22"
$ws.Range("N19").Value = "0.9 °C 5:10 TU"
$ws.Range("E20").Value = "2026-02-20 05:49:25"
$ws.Range("N20").Value = "-5.7 °C 5:06 TU"
$ws.Range("E21").Value = "2026-02-20 05:49:28"
$ws.Range("H21").Value = "'45%"
$ws.Range("J21").Value = "1020.0 hPa"
$ws.Range("O21").Value = "6.4 °C"
$ws.Range("E22").Value = "2026-02-20 05:49:30"
$ws.Range("H22").Value = "'64%"
$ws.Range("O22").Value = "-6.6 °C"
$ws.Range("E23").Value = "2026-02-20 05:49:32"
$ws.Range("G23").Value = "212 cm"
$ws.Range("I23").Value = "3.0 mm"
$ws.Range("E24").Value = "2026-02-20 05:49:35"
$ws.Range("J24").Value = "1023.1 hPa"
$ws.Range("N24").Value = "5.1 °C 5:00 TU"
$ws.Range("E25").Value = "2026-02-20 05:49:38"
$ws.Range("I25").Value = "4.3 mm"
$ws.Range("L25").Value = "58.0 km/h - 15º 5:06 TU"
$ws.Range("N25").Value = "-5.7 °C 5:16 TU"
$ws.Range("O25").Value = "-4.9 °C"
$ws.Range("E26").Value = "2026-02-20 05:49:41"
$ws.Range("G26").Value = "2 cm"
$ws.Range("H26").Value = "'42%"
$ws.Range("J26").Value = "1018.2 hPa"
$ws.Range("O26").Value = "3.4 °C"
$ws.Range("E27").Value = "2026-02-20 05:49:43"
$ws.Range("H27").Value = "'51%"
$ws.Range("O27").Value = "-2.9 °C"
$ws.Range("E28").Value = "2026-02-20 05:49:46"
$ws.Range("J28").Value = "1019.5 hPa"
$ws.Range("N28").Value = "0.7 °C 5:09 TU"
$ws.Range("O28").Value = "2.7 °C"
$ws.Range("E29").Value = "2026-02-20 05:49:49"
$ws.Range("N29").Value = "1.5 °C 5:12 TU"
$ws.Range("O29").Value = "3.5 °C"
$ws.Range("E30").Value = "2026-02-20 05:49:51"
$ws.Range("H30").Value = "'65%"
$ws.Range("J30").Value = "1018.0 hPa"
$ws.Range("N30").Value = "5.9 °C 5:29 TU"
$ws.Range("O30").Value = "9.0 °C"
$ws.Range("E31").Value = "2026-02-20 05:49:53"
$ws.Range("H31").Value = "'56%"
$ws.Range("J31").Value = "1016.5 hPa"
$ws.Range("K31").Value = "-0.1 MJ/m2"
$ws.Range("N31").Value = "9.6 °C 5:25 TU"
$ws.Range("O31").Value = "10.6 °C"
$ws.Range("E32").Value = "2026-02-20 05:49:56"
$ws.Range("H32").Value = "'89%"
$ws.Range("N32").Value = "1.6 °C 5:00 TU"
$ws.Range("E33").Value = "2026-02-20 05:49:59"
$ws.Range("J33").Value = "1019.0 hPa"
$ws.Range("O33").Value = "4.6 °C"
$ws.Range("E34").Value = "2026-02-20 05:50:01"
$ws.Range("H34").Value = "'64%"
$ws.Range("E35").Value = "2026-02-20 05:50:04"
$ws.Range("J35").Value = "1023.8 hPa"
$ws.Range("N35").Value = "1.7 °C 5:29 TU"
$ws.Range("E36").Value = "2026-02-20 05:50:06"
$ws.Range("J36").Value = "1018.6 hPa"
$ws.Range("L36").Value = "78.8 km/h - 325º 5:29 TU"
$ws.Range("E37").Value = "2026-02-20 05:50:09"
$ws.Range("H37").Value = "'67%"
$ws.Range("J37").Value = "1021.1 hPa"
$ws.Range("N37").Value = "-0.2 °C 5:29 TU"
$ws.Range("O37").Value = "2.0 °C"
$ws.Range("E38").Value = "2026-02-20 05:50:12"
$ws.Range("H38").Value = "'75%"
$ws.Range("N38").Value = "2.4 °C 5:29 TU"
$ws.Range("O38").Value = "4.3 °C"
$ws.Range("E39").Value = "2026-02-20 05:50:14"
$ws.Range("H39").Value = "'75%"
$ws.Range("E40").Value = "2026-02-20 05:50:17"
$ws.Range("H40").Value = "'49%"
$ws.Range("J40").Value = "1020.9 hPa"
$ws.Range("O40").Value = "7.2 °C"
$ws.Range("E41").Value = "2026-02-20 05:50:20"
$ws.Range("H41").Value = "'51%"
$ws.Range("N41").Value = "10.1 °C 5:07 TU"
$ws.Range("E42").Value = "2026-02-20 05:50:23"
$ws.Range("H42").Value = "'93%"
$ws.Range("N42").Value = "3.0 °C 5:24 TU"
$ws.Range("O42").Value = "4.6 °C"
$ws.Range("E43").Value = "2026-02-20 05:50:25"
$ws.Range("H43").Value = "'85%"
$ws.Range("N43").Value = "-0.1 °C 5:29 TU"
$ws.Range("O43").Value = "1.3 °C"
$ws.Range("E44").Value = "2026-02-20 05:50:27"
$ws.Range("I44").Value = "4.1 mm"
$ws.Range("O44").Value = "-6.0 °C"
$ws.Range("E45").Value = "2026-02-20 05:50:30"
$ws.Range("J45").Value = "1027.6 hPa"
$ws.Range("N45").Value = "1.7 °C 5:21 TU"
$ws.Range("E46").Value = "2026-02-20 05:50:33"
$ws.Range("J46").Value = "1023.8 hPa"
$ws.Range("N46").Value = "8.7 °C 5:28 TU"
